$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (and its "Through YYYY-MM-DD" label) to reflect the new date.
$ws.Name = "Through 2021-10-07"

# Update the October row label to reflect the new "through" date.
$ws.Range("A11").Value = "October (through 10-07)"

# Update the October (row 11) figures.
$ws.Range("C11").Value = 10
$ws.Range("E11").Value = 16
$ws.Range("G11").Value = 35
$ws.Range("H11").Value = 51

# Update the Total (row 12) figures.
$ws.Range("C12").Value = 439
$ws.Range("E12").Value = 564
$ws.Range("G12").Value = 936
$ws.Range("H12").Value = 1300
